$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $val) {
    $cell = $ws.Range($cellref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

$ws.Range('D2').Value = '23.947.91'
$ws.Range('E2').Value = '  +16.36%  '
$ws.Range('D3').Value = '1.661.97'
$ws.Range('E3').Value = '  +12.65%  '
Set-TextValue 'D4' '1.002'
$ws.Range('E4').Value = '  -0.73%  '
Set-TextValue 'D5' '308.25'
$ws.Range('E5').Value = '  +11.10%  '
Set-TextValue 'D6' '0.9938'
$ws.Range('E6').Value = '  +3.75%  '
Set-TextValue 'D7' '0.3716'
$ws.Range('E7').Value = '  +4.96%  '
Set-TextValue 'D8' '0.3445'
$ws.Range('E8').Value = '  +12.04%  '
Set-TextValue 'D9' '47.78'
$ws.Range('E9').Value = '  +21.48%  '
Set-TextValue 'D10' '1.170'
$ws.Range('E10').Value = '  +8.06%  '
Set-TextValue 'D11' '0.07246'
$ws.Range('E11').Value = '  +9.11%  '
Set-TextValue 'D12' '0.9950'
$ws.Range('E12').Value = '  -0.83%  '
Set-TextValue 'D13' '20.62'
$ws.Range('E13').Value = '  +14.03%  '
Set-TextValue 'D14' '6.033'
$ws.Range('E14').Value = '  +10.35%  '
Set-TextValue 'D15' '6.754'
$ws.Range('D16').Value = '1.660.90'
$ws.Range('E16').Value = '  +12.76%  '
Set-TextValue 'D17' '0.00001095'
$ws.Range('E17').Value = '  +7.85%  '
Set-TextValue 'D18' '0.9933'
$ws.Range('E18').Value = '  +3.57%  '
Set-TextValue 'D19' '0.06708'
$ws.Range('E19').Value = '  +12.17%  '
Set-TextValue 'D20' '81.38'
$ws.Range('E20').Value = '  +18.15%  '
Set-TextValue 'D21' '16.43'
$ws.Range('E21').Value = '  +13.48%  '
Set-TextValue 'D22' '6.103'
$ws.Range('E22').Value = '  +11.30%  '
Set-TextValue 'D23' '11.99'
$ws.Range('E23').Value = '  +7.62%  '
$ws.Range('D24').Value = '23.965.64'
$ws.Range('E24').Value = '  +16.41%  '
Set-TextValue 'D25' '2.372'
$ws.Range('E25').Value = '  +3.98%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D26' '2.694'
$ws.Range('E26').Value = '  +29.68%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D27' '3.376'
$ws.Range('E27').Value = '  -8.40%  '
$ws.Range('E28').Value = '  +3.40%  '
Set-TextValue 'D29' '19.55'
$ws.Range('E29').Value = '  +13.99%  '
$ws.Range('D30').Value = '1.841.61'
$ws.Range('E30').Value = '  +12.81%  '
Set-TextValue 'D31' '126.95'
$ws.Range('E31').Value = '  +10.92%  '
Set-TextValue 'D32' '4.139'
$ws.Range('E32').Value = '  +4.87%  '
Set-TextValue 'D33' '6.286'
$ws.Range('E33').Value = '  +27.63%  '
Set-TextValue 'D34' '0.9788'
$ws.Range('E34').Value = '  +23.18%  '
$ws.Range('E35').Value = '  +20.64%  '
Set-TextValue 'D36' '0.08388'
$ws.Range('E36').Value = '  +6.03%  '
Set-TextValue 'D37' '12.40'
$ws.Range('E37').Value = '  +20.92%  '
Set-TextValue 'D38' '8.987'
$ws.Range('E38').Value = '  +21.60%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D39' '0.06352'
$ws.Range('E39').Value = '  +11.88%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D40' '5.306'
$ws.Range('E40').Value = '  +13.10%  '
Set-TextValue 'D41' '1.286'
$ws.Range('E41').Value = '  +7.20%  '
Set-TextValue 'D42' '0.02318'
$ws.Range('E42').Value = '  +14.91%  '
Set-TextValue 'D43' '0.2075'
$ws.Range('E43').Value = '  +12.59%  '
Set-TextValue 'D44' '0.6112'
$ws.Range('E44').Value = '  +17.09%  '
Set-TextValue 'D45' '0.9933'
$ws.Range('E45').Value = '  +3.64%  '
Set-TextValue 'D46' '3.823'
$ws.Range('E46').Value = '  +8.78%  '
Set-TextValue 'D47' '13.26'
$ws.Range('E47').Value = '  +11.13%  '
Set-TextValue 'D48' '0.5946'
$ws.Range('E48').Value = '  +15.43%  '
Set-TextValue 'D49' '127.37'
$ws.Range('E49').Value = '  +6.09%  '
Set-TextValue 'D50' '2.002'
$ws.Range('E50').Value = '  +10.72%  '
Set-TextValue 'D51' '0.07073'
$ws.Range('E51').Value = '  +10.47%  '
